# 1D Lateral interactions implementation - first iteration
# Adds new CONST rows for quartz/organic density, and expands the Mualem Van
# Genuchten soil-hydraulic-parameter table (alpha/n/residual water content)
# with water, clay and peat columns plus a Darcy friction factor / air
# tortuosity block used by the Carman-Kozeny model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: insert two new rows for rho_m / rho_o right after the
#     existing rho_w / rho_i block (old row 9 blank separator shifts down,
#     and every row below shifts by +2 automatically). ---
$ws.Rows("9:10").Insert()

$ws.Range("A9").Value = "rho_m"
$ws.Range("A10").Value = "rho_o"
$ws.Range("D9").Value = "density quartz"
$ws.Range("D10").Value = "not sure"
$ws.Range("B9").Value = 2650
$ws.Range("B10").Value = 2650

# --- Replace the old Mualem Van Genuchten sand/silt-only block (previously
#     rows 36-41, now rows 38-43 after the shift above) with an expanded
#     table covering water, sand, silt, clay and peat (columns entered
#     grouped by soil type: peat, then water, then clay). ---

# peat column
$ws.Range("A42").Value = "alpha_peat"
$ws.Range("A47").Value = "n_peat"
$ws.Range("A52").Value = "residual_wc_peat"
$ws.Range("D42").Value = "from Hydraulic properties of fen peat soils in Poland, Gnatowski 2010"

# water column
$ws.Range("A38").Value = "alpha_water"
$ws.Range("A43").Value = "n_water"
$ws.Range("A48").Value = "residual_wc_water"

# clay column
$ws.Range("A41").Value = "alpha_clay"
$ws.Range("A46").Value = "n_clay"
$ws.Range("A51").Value = "residual_wc_clay"

# sand / silt columns keep their existing labels (already in sharedStrings)
$ws.Range("A39").Value = "alpha_sand"
$ws.Range("A40").Value = "alpha_silt"
$ws.Range("A44").Value = "n_sand"
$ws.Range("A45").Value = "n_silt"
$ws.Range("A49").Value = "residual_wc_sand"
$ws.Range("A50").Value = "residual_wc_silt"

# values: alpha [1/m]
$ws.Range("B38").Value = 400
$ws.Range("B39").Value = 4.06
$ws.Range("B40").Value = 0.65
$ws.Range("B41").Value = 1.49
$ws.Range("B42").Value = 2.31

# values: n [-]
$ws.Range("B43").Value = 2.5
$ws.Range("B44").Value = 2
$ws.Range("B45").Value = 1.7
$ws.Range("B46").Value = 1.25
$ws.Range("B47").Value = 1.292

# values: residual water content [-]
$ws.Range("B48").Value = 0
$ws.Range("B49").Value = 0
$ws.Range("B50").Value = 0
$ws.Range("B51").Value = 0
$ws.Range("B52").Value = 0

# Apply the "0.00" number format to the whole new alpha/n/residual_wc block.
$ws.Range("B38:B52").NumberFormat = "0.00"

# --- New Darcy friction factor / air tortuosity block (separated by a
#     blank row 53). ---
$ws.Range("A54").Value = "Darcy_friction_factor"
$ws.Range("C54").Value = "rough-pipe regime"
$ws.Range("B54").Value = 0.1

$ws.Range("A55").Value = "tortuosity_air"
$ws.Range("C55").Value = "used in Carman Kozeny model"
$ws.Range("B55").Value = 2.5

# Match the author's final selection / view state.
$ws.Range("A54:D55").Select()
